$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.283.65"
$ws.Range("E2").Value = "  -1.05%  "

$ws.Range("D3").Value = "3.079.42"
$ws.Range("E3").Value = "  -1.29%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.73%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.078.70"
$ws.Range("E8").Value = "  -1.18%  "

$ws.Range("E9").Value = "  -2.33%  "

$ws.Range("E10").Value = "  -1.74%  "

$ws.Range("E11").Value = "  -2.63%  "

$ws.Range("E12").Value = "  -2.76%  "

$ws.Range("E13").Value = "  -3.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.55%  "

$ws.Range("E15").Value = "  -1.44%  "

$ws.Range("D16").Value = "3.593.98"
$ws.Range("E16").Value = "  -1.15%  "

$ws.Range("D17").Value = "66.214.75"
$ws.Range("E17").Value = "  -1.05%  "

$ws.Range("E18").Value = "  -3.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.34%  "

$ws.Range("D20").Value = "3.079.20"
$ws.Range("E20").Value = "  -1.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "487.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.686"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.13%  "

$ws.Range("E26").Value = "  -3.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.26%  "

$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.36%  "

$ws.Range("E30").Value = "  -4.94%  "

$ws.Range("E31").Value = "  -2.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.89%  "

$ws.Range("E33").Value = "  -3.52%  "

$ws.Range("D34").Value = "0.0₃0902"
$ws.Range("E34").Value = "  -4.65%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("B36").Value = "Arweave"
$ws.Range("C36").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "47.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.39%  "

$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.948"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.90%  "

$ws.Range("E39").Value = "  -0.63%  "

$ws.Range("E40").Value = "  -4.65%  "

$ws.Range("E41").Value = "  -4.37%  "

$ws.Range("E42").Value = "  -4.74%  "

$ws.Range("D43").Value = "2.783.86"
$ws.Range("E43").Value = "  -1.19%  "

$ws.Range("E44").Value = "  -2.83%  "

$ws.Range("E45").Value = "  -1.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "365.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.85%  "

$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.74%  "

$ws.Range("E50").Value = "  -2.06%  "

$ws.Range("E51").Value = "  -2.09%  "

